$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from 45175 (2023-09-06) to 45183 (2023-09-14)
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45183
}
